$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 42.409254
$ws.Range("H2").Value = 127.227762
$ws.Range("I2").Value = 0.6138221220752584
$ws.Range("J2").Value = 0.6138221220752584
$ws.Range("M2").Value = 1.443038
$ws.Range("N2").Value = 4.329114
$ws.Range("O2").Value = 0.0289666880885598
$ws.Range("P2").Value = 0.0289666880885598
$ws.Range("Q2").Value = 61.198165073652
$ws.Range("R2").Value = 550.783485662868
$ws.Range("S2").Value = 0.01778039395201189
$ws.Range("T2").Value = 0.01778039395201188
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 42.409254
$ws.Range("H3").Value = 127.227762
$ws.Range("I3").Value = 0.6138221220752584
$ws.Range("J3").Value = 0.6138221220752584
$ws.Range("N3").Value = 87.61054300000001
$ws.Range("O3").Value = 0.5862140087672342
$ws.Range("P3").Value = 0.5862140087672342
$ws.Range("Q3").Value = 1238.499257054974
$ws.Range("R3").Value = 11146.49331349477
$ws.Range("S3").Value = 0.3598311268517478
$ws.Range("T3").Value = 0.3598311268517478
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 42.409254
$ws.Range("H4").Value = 127.227762
$ws.Range("I4").Value = 0.6138221220752584
$ws.Range("J4").Value = 0.6138221220752584
$ws.Range("M4").Value = 19.170603
$ws.Range("N4").Value = 57.511809
$ws.Range("O4").Value = 0.384819303144206
$ws.Range("P4").Value = 0.384819303144206
$ws.Range("Q4").Value = 813.0109719601621
$ws.Range("R4").Value = 7317.098747641458
$ws.Range("S4").Value = 0.2362106012714987
$ws.Range("T4").Value = 0.2362106012714987
$ws.Range("I5").Value = 0.07014398987036251
$ws.Range("J5").Value = 0.07014398987036251
$ws.Range("M5").Value = 1.443038
$ws.Range("N5").Value = 4.329114
$ws.Range("O5").Value = 0.0289666880885598
$ws.Range("P5").Value = 0.0289666880885598
$ws.Range("Q5").Value = 6.993367160665332
$ws.Range("R5").Value = 62.94030444598799
$ws.Range("S5").Value = 0.002031839075861889
$ws.Range("T5").Value = 0.002031839075861889
$ws.Range("I6").Value = 0.07014398987036251
$ws.Range("J6").Value = 0.07014398987036251
$ws.Range("N6").Value = 87.61054300000001
$ws.Range("O6").Value = 0.5862140087672342
$ws.Range("P6").Value = 0.5862140087672342
$ws.Range("S6").Value = 0.04111938949283347
$ws.Range("T6").Value = 0.04111938949283347
$ws.Range("I7").Value = 0.07014398987036251
$ws.Range("J7").Value = 0.07014398987036251
$ws.Range("M7").Value = 19.170603
$ws.Range("N7").Value = 57.511809
$ws.Range("O7").Value = 0.384819303144206
$ws.Range("P7").Value = 0.384819303144206
$ws.Range("Q7").Value = 92.90612268724199
$ws.Range("R7").Value = 836.1551041851779
$ws.Range("S7").Value = 0.02699276130166715
$ws.Range("T7").Value = 0.02699276130166714
$ws.Range("G8").Value = 21.83492733333334
$ws.Range("H8").Value = 65.50478200000001
$ws.Range("I8").Value = 0.3160338880543792
$ws.Range("J8").Value = 0.3160338880543791
$ws.Range("M8").Value = 1.443038
$ws.Range("N8").Value = 4.329114
$ws.Range("O8").Value = 0.0289666880885598
$ws.Range("P8").Value = 0.0289666880885598
$ws.Range("Q8").Value = 31.50862986923867
$ws.Range("R8").Value = 283.577668823148
$ws.Range("S8").Value = 0.009154455060686026
$ws.Range("T8").Value = 0.009154455060686023
$ws.Range("G9").Value = 21.83492733333334
$ws.Range("H9").Value = 65.50478200000001
$ws.Range("I9").Value = 0.3160338880543792
$ws.Range("J9").Value = 0.3160338880543791
$ws.Range("N9").Value = 87.61054300000001
$ws.Range("O9").Value = 0.5862140087672342
$ws.Range("P9").Value = 0.5862140087672342
$ws.Range("Q9").Value = 637.6566133462919
$ws.Range("R9").Value = 5738.909520116627
$ws.Range("S9").Value = 0.1852634924226529
$ws.Range("T9").Value = 0.1852634924226529
$ws.Range("G10").Value = 21.83492733333334
$ws.Range("H10").Value = 65.50478200000001
$ws.Range("I10").Value = 0.3160338880543792
$ws.Range("J10").Value = 0.3160338880543791
$ws.Range("M10").Value = 19.170603
$ws.Range("N10").Value = 57.511809
$ws.Range("O10").Value = 0.384819303144206
$ws.Range("P10").Value = 0.384819303144206
$ws.Range("Q10").Value = 418.5887234411821
$ws.Range("R10").Value = 3767.298510970638
$ws.Range("S10").Value = 0.1216159405710402
$ws.Range("T10").Value = 0.1216159405710402
